$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Agosto de 2020 a las 07:11"

# Uzbekistan's new case counts overtake Venezuela's (which are unchanged),
# so the two countries swap places in the ranking (row 61 <-> row 62).
$ws.Range("A61").Value = "Uzbekistan"
$ws.Range("B61").Value = 37693
$ws.Range("C61").Value = 146
$ws.Range("D61").Value = 33442
$ws.Range("E61").Value = 3997
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 254

$ws.Range("A62").Value = "Venezuela"
$ws.Range("B62").Value = 37567
$ws.Range("D62").Value = 26330
$ws.Range("E62").Value = 10926
$ws.Range("H62").Value = 311

# India (row 6)
$ws.Range("B6").Value = 2905823
$ws.Range("C6").Value = 1494
$ws.Range("D6").Value = 2158946
$ws.Range("E6").Value = 691902

# Pakistan (row 18)
$ws.Range("B18").Value = 291588
$ws.Range("C18").Value = 630
$ws.Range("D18").Value = 273579
$ws.Range("E18").Value = 11790
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 6219

# Tailandia (row 118)
$ws.Range("B118").Value = 3390
$ws.Range("C118").Value = 1
$ws.Range("D118").Value = 3219

# Islas Turcas y Caicos (row 180)
$ws.Range("B180").Value = 334
$ws.Range("C180").Value = 7
$ws.Range("D180").Value = 102
$ws.Range("E180").Value = 230

# Mongolia (row 182)
$ws.Range("D182").Value = 283
$ws.Range("E182").Value = 15

# Butan (row 189)
$ws.Range("B189").Value = 153
$ws.Range("C189").Value = 3
$ws.Range("D189").Value = 108
